# Add a "Save" column (H) to the s_vals sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header - copy formatting (style) from the neighboring "sum" header cell
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Save flags per row (2-34): 1 for the top-3 "sum" rows (13, 19, 34), else 0
$saveValues = @{
    2  = 0
    3  = 0
    4  = 0
    5  = 0
    6  = 0
    7  = 0
    8  = 0
    9  = 0
    10 = 0
    11 = 0
    12 = 0
    13 = 1
    14 = 0
    15 = 0
    16 = 0
    17 = 0
    18 = 0
    19 = 1
    20 = 0
    21 = 0
    22 = 0
    23 = 0
    24 = 0
    25 = 0
    26 = 0
    27 = 0
    28 = 0
    29 = 0
    30 = 0
    31 = 0
    32 = 0
    33 = 0
    34 = 1
}

foreach ($row in $saveValues.Keys) {
    $ws.Cells.Item($row, 8).Value = $saveValues[$row]
}
